# "Generate Report for Handback" - populate handback info for zh-cn / de-de
# localization rows, update the overview status text, and widen the
# columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$mdFile   = "14c5eaed-5d6a-401f-8390-6458a266a54d.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6aca54f5835fc263c57e2bc234551a42320d7b0/e2e/14c5eaed-5d6a-401f-8390-6458a266a54d.md"
$zhXlf    = "14c5eaed-5d6a-401f-8390-6458a266a54d.e700243b4f973f6f0fba52dd76dfe2e7c8fe8655.zh-cn.xlf"
$deXlf    = "14c5eaed-5d6a-401f-8390-6458a266a54d.e700243b4f973f6f0fba52dd76dfe2e7c8fe8655.de-de.xlf"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (Overview E2/F2 and the Status column on each language sheet)
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws2.Range("C2").Value = $newStatus
$ws3.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn row 2: Latest Target File / Latest Handback File / Latest Handback DateTime
# ---------------------------------------------------------------------
$ws2.Range("J2").Value = $zhXlf
$ws2.Range("K2").Value = "2016-11-09 10:42:43"
$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdFile)

# ---------------------------------------------------------------------
# de-de row 2: Latest Target File / Latest Handback File / Latest Handback DateTime
# ---------------------------------------------------------------------
$ws3.Range("J2").Value = $deXlf
$ws3.Range("K2").Value = "2016-11-09 10:43:02"
$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdFile)

# ---------------------------------------------------------------------
# Column widths: widen the columns now holding the longer handback text
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.2   # E - zh-cn status
$ws1.Columns.Item(6).ColumnWidth = 29.2   # F - de-de status

foreach ($ws in @($ws2, $ws3)) {
    $ws.Columns.Item(3).ColumnWidth = 29.2    # C - Status
    $ws.Columns.Item(9).ColumnWidth = 39.2    # I - Latest Target File
    $ws.Columns.Item(10).ColumnWidth = 39.2   # J - Latest Handback File
}
